$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 973546.9
$ws.Range("I70").Value = 3889338
$ws.Range("K70").Value = 11668014
$ws.Range("M70").Value = -11667744
$ws.Range("H73").Value = 973546.9
$ws.Range("I73").Value = 3889338
$ws.Range("K73").Value = 11668014
$ws.Range("M73").Value = -11667078
$ws.Range("H103").Value = 964.06665
$ws.Range("I103").Value = 326.33334
$ws.Range("J103").Value = 1123.5
$ws.Range("K103").Value = 979.0000200000001
$ws.Range("L103").Value = 3370.5
$ws.Range("M103").Value = -393.0000200000001
$ws.Range("N103").Value = -4542.5
$ws.Range("H107").Value = 387.42856
$ws.Range("I107").Value = 387.42856
$ws.Range("K107").Value = 387.42856
$ws.Range("M107").Value = 1532.57144
$ws.Range("H137").Value = 2352.5
$ws.Range("I137").Value = 2169.5
$ws.Range("K137").Value = 6508.5
$ws.Range("M137").Value = -3958.5
$ws.Range("H138").Value = 2641.16
$ws.Range("I138").Value = 1573.5
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 4720.5
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 419.5
$ws.Range("N138").Value = -22280

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1404.6428
$ws.Range("I2").Value = 1394.4546
$ws.Range("J2").Value = 1442
$ws.Range("K2").Value = 1394.4546
$ws.Range("L2").Value = 1442
$ws.Range("M2").Value = -1281.4546
$ws.Range("N2").Value = -1668
$ws.Range("H61").Value = 15152944
$ws.Range("I61").Value = 16667918
$ws.Range("K61").Value = 16667918
$ws.Range("M61").Value = -16667706
$ws.Range("H110").Value = 76947400
$ws.Range("I110").Value = 90910660
$ws.Range("K110").Value = 90910660
$ws.Range("M110").Value = -90908615
$ws.Range("H116").Value = 1404.6428
$ws.Range("I116").Value = 1394.4546
$ws.Range("J116").Value = 1442
$ws.Range("K116").Value = 1394.4546
$ws.Range("L116").Value = 1442
$ws.Range("M116").Value = 899.5454
$ws.Range("N116").Value = -6030
$ws.Range("H136").Value = 15152944
$ws.Range("I136").Value = 16667918
$ws.Range("K136").Value = 50003754
$ws.Range("M136").Value = -50001204

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1404.6428
$ws.Range("I3").Value = 1394.4546
$ws.Range("J3").Value = 1442
$ws.Range("K3").Value = 1394.4546
$ws.Range("L3").Value = 1442
$ws.Range("M3").Value = -1280.4546
$ws.Range("N3").Value = -1670

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1378.5385
$ws.Range("I16").Value = 1123.6875
$ws.Range("K16").Value = 1123.6875
$ws.Range("M16").Value = -836.6875
$ws.Range("H22").Value = 310.2
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 325.25
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 325.25
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -1025.25
$ws.Range("H39").Value = 2975
$ws.Range("I39").Value = 2975
$ws.Range("K39").Value = 2975
$ws.Range("M39").Value = -2584
$ws.Range("H49").Value = 2975
$ws.Range("I49").Value = 2975
$ws.Range("K49").Value = 2975
$ws.Range("M49").Value = -2793
$ws.Range("H86").Value = 60236.2
$ws.Range("I86").Value = 96871.55499999999
$ws.Range("J86").Value = 5283.1665
$ws.Range("K86").Value = 96871.55499999999
$ws.Range("L86").Value = 5283.1665
$ws.Range("M86").Value = -95748.55499999999
$ws.Range("N86").Value = -7529.1665
$ws.Range("H89").Value = 60236.2
$ws.Range("I89").Value = 96871.55499999999
$ws.Range("J89").Value = 5283.1665
$ws.Range("K89").Value = 484357.775
$ws.Range("L89").Value = 26415.8325
$ws.Range("M89").Value = -478741.775
$ws.Range("N89").Value = -37647.8325
$ws.Range("H107").Value = 2121.7058
$ws.Range("J107").Value = 548.3333
$ws.Range("L107").Value = 548.3333
$ws.Range("N107").Value = -4388.3333
$ws.Range("H113").Value = 1378.5385
$ws.Range("I113").Value = 1123.6875
$ws.Range("K113").Value = 1123.6875
$ws.Range("M113").Value = 1046.3125
$ws.Range("H132").Value = 2170.7
$ws.Range("I132").Value = 1813.75
$ws.Range("K132").Value = 5441.25
$ws.Range("M132").Value = -2911.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17811550
$ws.Range("I4").Value = 45116324
$ws.Range("J4").Value = 4087.1304
$ws.Range("K4").Value = 135348972
$ws.Range("L4").Value = 12261.3912
$ws.Range("M4").Value = -135348860
$ws.Range("N4").Value = -12485.3912
$ws.Range("H86").Value = 2257.1428
$ws.Range("I86").Value = 850.5
$ws.Range("K86").Value = 2551.5
$ws.Range("M86").Value = -1365.5
$ws.Range("H89").Value = 2257.1428
$ws.Range("I89").Value = 850.5
$ws.Range("K89").Value = 7654.5
$ws.Range("M89").Value = -1726.5
$ws.Range("H122").Value = 566.3333
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 68500
$ws.Range("I62").Value = 47500
$ws.Range("K62").Value = 47500
$ws.Range("M62").Value = -46814
$ws.Range("H65").Value = 68500
$ws.Range("I65").Value = 47500
$ws.Range("K65").Value = 142500
$ws.Range("M65").Value = -139068
$ws.Range("H113").Value = 2900.9285
$ws.Range("I113").Value = 980
$ws.Range("J113").Value = 3968.111
$ws.Range("K113").Value = 980
$ws.Range("L113").Value = 3968.111
$ws.Range("M113").Value = 1190
$ws.Range("N113").Value = -8308.111000000001
$ws.Range("H132").Value = 3416.7222
$ws.Range("I132").Value = 2983.0344
$ws.Range("K132").Value = 8949.1032
$ws.Range("M132").Value = -6419.1032
$ws.Range("H135").Value = 99933.336
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 50002880
$ws.Range("J7").Value = 4914.143
$ws.Range("K7").Value = 50002880
$ws.Range("L7").Value = 4914.143
$ws.Range("M7").Value = -50002768
$ws.Range("N7").Value = -5138.143
$ws.Range("H100").Value = 2331.1304
$ws.Range("I100").Value = 2086.077
$ws.Range("J100").Value = 2649.7
$ws.Range("K100").Value = 2086.077
$ws.Range("L100").Value = 2649.7
$ws.Range("M100").Value = -1545.077
$ws.Range("N100").Value = -3731.7
$ws.Range("I126").Value = 50002880
$ws.Range("J126").Value = 4914.143
$ws.Range("K126").Value = 150008640
$ws.Range("L126").Value = 14742.429
$ws.Range("M126").Value = -150006170
$ws.Range("N126").Value = -19682.429

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 591.5333000000001
$ws.Range("I107").Value = 598
$ws.Range("J107").Value = 549.5
$ws.Range("K107").Value = 1794
$ws.Range("L107").Value = 1648.5
$ws.Range("M107").Value = 126
$ws.Range("N107").Value = -5488.5
$ws.Range("H136").Value = 6687.875
$ws.Range("I136").Value = 2252
$ws.Range("J136").Value = 8166.5
$ws.Range("K136").Value = 6756
$ws.Range("L136").Value = 24499.5
$ws.Range("M136").Value = -4206
$ws.Range("N136").Value = -29599.5

# --- Clear removed cells (values dropped entirely in the diff) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M135").ClearContents()
